$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 258, shifting existing rows 258:356 down to 259:357
$ws.Rows(258).Insert()

# Populate the newly inserted row 258 with the new record's data
$ws.Cells.Item(258, 1).Value = 6
$ws.Cells.Item(258, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(258, 3).Value = "Metropolitana"
$ws.Cells.Item(258, 4).Value = 44627
$ws.Cells.Item(258, 5).Value = 13
$ws.Cells.Item(258, 6).Value = 100112032
$ws.Cells.Item(258, 7).Value = "Zapallo italiano"
$ws.Cells.Item(258, 8).Value = "Sin especificar"
$ws.Cells.Item(258, 9).Value = "Primera"
$ws.Cells.Item(258, 10).Value = 930
$ws.Cells.Item(258, 11).Value = 11000
$ws.Cells.Item(258, 12).Value = 12000
$ws.Cells.Item(258, 13).Value = 11602
$ws.Cells.Item(258, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(258, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(258, 16).Value = 232
$ws.Cells.Item(258, 17).Value = 50
$ws.Cells.Item(258, 18).Value = "Hortaliza"
